$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the
# Heading1 paragraph ("Play Atlantean GigaRise for Free: Read Our Review").
# The new paragraph has a leading empty run, a bold "Meta description" run,
# and a normal run with the rest of the meta-description text.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Discover the features of Atlantean GigaRise, a highly volatile slot game with up to 294,912 ways to win. Play for free and read our review to learn more.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$titlePara.Range.InsertXML($metaXml, 1)

# ---------------------------------------------------------------------------
# Change 2: remove the duplicated bold "Play Atlantean GigaRise for Free:
# Read Our Review" paragraph that used to sit right before the final
# (italic) meta-description paragraph at the end of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupeHeadingPara = $d.Paragraphs.Item($count - 1)
$dupeHeadingPara.Range.Delete()

# ---------------------------------------------------------------------------
# Change 3: replace the text of the final (italic) paragraph -- it used to
# hold the meta description, and now should hold the image-generation
# prompt instead. We only replace the text up to (not including) the
# paragraph mark so the run's formatting (italic) and the leading empty run
# are preserved.
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count2)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalRange.Text = 'Please create a feature image for "Atlantean GigaRise" that features a happy Maya warrior with glasses in a cartoon style. The image should depict the warrior standing next to the underwater scene of Atlantis with the mountains and colonnades in the background. The warrior should be holding a trident with a smile on their face and wearing glasses. The image should be colorful and eye-catching to attract the attention of players.'
